$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after row 96 (this pushes the old row 97 and everything
# below it down by one row, matching the target diff where a brand new record is
# inserted before the former row 97, and the very last former row (127) ends up
# at row 128).
$ws.Rows.Item(97).Insert()

# Copy the static / shared columns from the row above (row 96) into the new row 97,
# since this dataset shares the same market/product metadata for every row.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(97, $col).Value = $ws.Cells.Item(96, $col).Value2
}

# Now set the specific values that differ for this new record.
$ws.Cells.Item(97, 4).Value = 44876       # D97 Fecha
$ws.Cells.Item(97, 11).Value = "Murcott"  # K97 Variedad
$ws.Cells.Item(97, 12).Value = "Segunda"  # L97 Calidad
$ws.Cells.Item(97, 13).Value = 270        # M97 Volumen
$ws.Cells.Item(97, 14).Value = 15000      # N97 Precio minimo
$ws.Cells.Item(97, 15).Value = 16000      # O97 Precio maximo
$ws.Cells.Item(97, 16).Value = 15500      # P97 Precio promedio ponderado
$ws.Cells.Item(97, 19).Value = 775        # S97 Precio $/Kg
$ws.Cells.Item(97, 20).Value = 20         # T97 Kg / unidad

# Apply the date style (style index 2, number format "YYYY-MM-DD HH:MM:SS") used
# by the rest of column D to the new D97 cell, matching the other date cells.
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(96, 4).NumberFormat
